# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled run).
# Price (col D) and Volume(1h) (col E) are stored as literal text, not
# numbers, in this sheet - so for any new price that happens to look like a
# plain number (e.g. "1.001"), the cell's NumberFormat is set to Text ("@")
# before the write; otherwise Excel COM auto-coerces the string to a real
# number on assignment. Percent strings in col E always contain surrounding
# spaces/"%" so they never round-trip as numbers and need no such guard.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.016.16'
$ws.Range('E2').Value = '  -1.55%  '

$ws.Range('D3').Value = '1.766.09'
$ws.Range('E3').Value = '  -3.51%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.48%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.53'
$ws.Range('E5').Value = '  -2.21%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.73%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4265'
$ws.Range('E7').Value = '  -4.59%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3617'
$ws.Range('E8').Value = '  -4.62%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.63'
$ws.Range('E9').Value = '  -3.83%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07448'
$ws.Range('E10').Value = '  -4.36%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.099'
$ws.Range('E11').Value = '  -3.75%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9996'
$ws.Range('E12').Value = '  +0.62%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.14'
$ws.Range('E13').Value = '  -5.32%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.069'
$ws.Range('E14').Value = '  -4.23%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.326'
$ws.Range('E15').Value = '  -3.02%  '

$ws.Range('D16').Value = '1.780.34'
$ws.Range('E16').Value = '  -2.54%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.18'
$ws.Range('E17').Value = '  -1.39%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001058'
$ws.Range('E18').Value = '  -2.53%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06383'
$ws.Range('E19').Value = '  +0.10%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9995'
$ws.Range('E20').Value = '  +0.49%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.05'
$ws.Range('E21').Value = '  -3.15%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.985'
$ws.Range('E22').Value = '  -6.03%  '

$ws.Range('D23').Value = '28.016.27'
$ws.Range('E23').Value = '  -1.73%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.26'
$ws.Range('E24').Value = '  -4.93%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.132'
$ws.Range('E25').Value = '  -1.77%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.51'
$ws.Range('E26').Value = '  +2.10%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.19'
$ws.Range('E27').Value = '  -3.89%  '

$ws.Range('D28').Value = '1.987.30'
$ws.Range('E28').Value = '  -2.38%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.135'
$ws.Range('E29').Value = '  -10.35%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.74'
$ws.Range('E30').Value = '  -4.07%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.156'
$ws.Range('E31').Value = '  -5.85%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.631'
$ws.Range('E32').Value = '  -4.18%  '

$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.602'
$ws.Range('E33').Value = '  -1.82%  '

$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08905'
$ws.Range('E34').Value = '  -3.92%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02315'
$ws.Range('E36').Value = '  -2.04%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2104'
$ws.Range('E37').Value = '  -4.43%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.025'
$ws.Range('E38').Value = '  -3.60%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06033'
$ws.Range('E39').Value = '  -3.72%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6373'
$ws.Range('E40').Value = '  -4.38%  '

$ws.Range('E41').Value = '  -0.49%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9993'
$ws.Range('E42').Value = '  +0.70%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.847'
$ws.Range('E43').Value = '  -3.13%  '

$ws.Range('B44').Value = 'WEMIXTOKEN'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.396'
$ws.Range('E44').Value = '  -0.93%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.37'
$ws.Range('E45').Value = '  -4.74%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5918'
$ws.Range('E46').Value = '  -3.43%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.686'
$ws.Range('E47').Value = '  -1.90%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.009'
$ws.Range('E48').Value = '  -1.82%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.53'
$ws.Range('E49').Value = '  -4.02%  '

$ws.Range('E50').Value = '  +3.01%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06866'
$ws.Range('E51').Value = '  -2.14%  '
